$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 576-595 with revised cumulative figures ---

$ws.Range("F576").Value = 28929

$ws.Range("F579").Value = 32604

$ws.Range("F580").Value = 28835

$ws.Range("F581").Value = 26987
$ws.Range("G581").Value = 479

$ws.Range("F582").Value = 25837

$ws.Range("F583").Value = 29217

$ws.Range("F584").Value = 13242
$ws.Range("G584").Value = 249

$ws.Range("F585").Value = 14909

$ws.Range("F586").Value = 33544
$ws.Range("G586").Value = 703

$ws.Range("F587").Value = 28177
$ws.Range("G587").Value = 554

$ws.Range("F588").Value = 25350

$ws.Range("F589").Value = 25349
$ws.Range("G589").Value = 469

$ws.Range("F590").Value = 29169
$ws.Range("G590").Value = 578

$ws.Range("F591").Value = 14452
$ws.Range("G591").Value = 415

$ws.Range("F592").Value = 18287
$ws.Range("G592").Value = 652

$ws.Range("F593").Value = 36882
$ws.Range("G593").Value = 1186

$ws.Range("F594").Value = 29670
$ws.Range("G594").Value = 816

$ws.Range("F595").Value = 27200
$ws.Range("G595").Value = 900

# --- Append new rows 596-599 for the latest daily stats ---

$newRows = @(
    @{ Row = 596; A = 44490; B = 449775; C = 14282; D = 3470; E = 12886; F = 28305; G = 926 },
    @{ Row = 597; A = 44491; B = 453231; C = 15276; D = 3456; E = 12895; F = 25975; G = 859 },
    @{ Row = 598; A = 44492; B = 456438; C = 12572; D = 3207; E = 12903; F = 12155; G = 567 },
    @{ Row = 599; A = 44493; B = 457431; C = 3854;  D = 993;  E = 12917; F = 10627; G = 538 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 1).NumberFormat = "yyyy-mm-dd"
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
}
